$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 190.664594
$ws.Range("H2").Value = 571.993782
$ws.Range("I2").Value = 0.2001939625490346
$ws.Range("J2").Value = 0.2001939625490346
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.3542
$ws.Range("N2").Value = 49.0626
$ws.Range("O2").Value = 0.3510578481048182
$ws.Range("P2").Value = 0.3510578481048182
$ws.Range("Q2").Value = 3118.1669031948
$ws.Range("R2").Value = 28063.5021287532
$ws.Range("S2").Value = 0.07027966169604066
$ws.Range("T2").Value = 0.07027966169604065
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 190.664594
$ws.Range("H3").Value = 571.993782
$ws.Range("I3").Value = 0.2001939625490346
$ws.Range("J3").Value = 0.2001939625490346
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.477892333333333
$ws.Range("N3").Value = 7.433676999999999
$ws.Range("O3").Value = 0.05319022332950721
$ws.Range("P3").Value = 0.0531902233295072
$ws.Range("Q3").Value = 472.4463357107126
$ws.Range("R3").Value = 4252.017021396414
$ws.Range("S3").Value = 0.01064836157720215
$ws.Range("T3").Value = 0.01064836157720215
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 190.664594
$ws.Range("H4").Value = 571.993782
$ws.Range("I4").Value = 0.2001939625490346
$ws.Range("J4").Value = 0.2001939625490346
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.62789
$ws.Range("N4").Value = 61.88367
$ws.Range("O4").Value = 0.4427965094191643
$ws.Range("P4").Value = 0.4427965094191643
$ws.Range("Q4").Value = 3933.008271926659
$ws.Range("R4").Value = 35397.07444733994
$ws.Range("S4").Value = 0.08864518782350345
$ws.Range("T4").Value = 0.08864518782350343
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 190.664594
$ws.Range("H5").Value = 571.993782
$ws.Range("I5").Value = 0.2001939625490346
$ws.Range("J5").Value = 0.2001939625490346
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.125502333333333
$ws.Range("N5").Value = 21.376507
$ws.Range("O5").Value = 0.1529554191465104
$ws.Range("P5").Value = 0.1529554191465104
$ws.Range("Q5").Value = 1358.581009431053
$ws.Range("R5").Value = 12227.22908487947
$ws.Range("S5").Value = 0.0306207514522884
$ws.Range("T5").Value = 0.0306207514522884
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 531.1103823333333
$ws.Range("H6").Value = 1593.331147
$ws.Range("I6").Value = 0.5576551459273178
$ws.Range("J6").Value = 0.5576551459273177
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 16.3542
$ws.Range("N6").Value = 49.0626
$ws.Range("O6").Value = 0.3510578481048182
$ws.Range("P6").Value = 0.3510578481048182
$ws.Range("Q6").Value = 8685.885414755801
$ws.Range("R6").Value = 78172.9687328022
$ws.Range("S6").Value = 0.1957692155138226
$ws.Range("T6").Value = 0.1957692155138225
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 531.1103823333333
$ws.Range("H7").Value = 1593.331147
$ws.Range("I7").Value = 0.5576551459273178
$ws.Range("J7").Value = 0.5576551459273177
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.477892333333333
$ws.Range("N7").Value = 7.433676999999999
$ws.Range("O7").Value = 0.05319022332950721
$ws.Range("P7").Value = 0.0531902233295072
$ws.Range("Q7").Value = 1316.034344537502
$ws.Range("R7").Value = 11844.30910083752
$ws.Range("S7").Value = 0.02966180175272297
$ws.Range("T7").Value = 0.02966180175272295
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 531.1103823333333
$ws.Range("H8").Value = 1593.331147
$ws.Range("I8").Value = 0.5576551459273178
$ws.Range("J8").Value = 0.5576551459273177
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 20.62789
$ws.Range("N8").Value = 61.88367
$ws.Range("O8").Value = 0.4427965094191643
$ws.Range("P8").Value = 0.4427965094191643
$ws.Range("Q8").Value = 10955.68654462994
$ws.Range("R8").Value = 98601.17890166947
$ws.Range("S8").Value = 0.246927752076251
$ws.Range("T8").Value = 0.246927752076251
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 531.1103823333333
$ws.Range("H9").Value = 1593.331147
$ws.Range("I9").Value = 0.5576551459273178
$ws.Range("J9").Value = 0.5576551459273177
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.125502333333333
$ws.Range("N9").Value = 21.376507
$ws.Range("O9").Value = 0.1529554191465104
$ws.Range("P9").Value = 0.1529554191465104
$ws.Range("Q9").Value = 3784.428268573725
$ws.Range("R9").Value = 34059.85441716353
$ws.Range("S9").Value = 0.08529637658452133
$ws.Range("T9").Value = 0.08529637658452131
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 79.06597733333334
$ws.Range("H10").Value = 237.197932
$ws.Range("I10").Value = 0.08301767503395074
$ws.Range("J10").Value = 0.08301767503395074
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.3542
$ws.Range("N10").Value = 49.0626
$ws.Range("O10").Value = 0.3510578481048182
$ws.Range("P10").Value = 0.3510578481048182
$ws.Range("Q10").Value = 1293.0608065048
$ws.Range("R10").Value = 11637.5472585432
$ws.Range("S10").Value = 0.02914400635208384
$ws.Range("T10").Value = 0.02914400635208383
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 79.06597733333334
$ws.Range("H11").Value = 237.197932
$ws.Range("I11").Value = 0.08301767503395074
$ws.Range("J11").Value = 0.08301767503395074
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.477892333333333
$ws.Range("N11").Value = 7.433676999999999
$ws.Range("O11").Value = 0.05319022332950721
$ws.Range("P11").Value = 0.0531902233295072
$ws.Range("Q11").Value = 195.9169790617738
$ws.Range("R11").Value = 1763.252811555964
$ws.Range("S11").Value = 0.004415728675352295
$ws.Range("T11").Value = 0.004415728675352295
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 79.06597733333334
$ws.Range("H12").Value = 237.197932
$ws.Range("I12").Value = 0.08301767503395074
$ws.Range("J12").Value = 0.08301767503395074
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 20.62789
$ws.Range("N12").Value = 61.88367
$ws.Range("O12").Value = 0.4427965094191643
$ws.Range("P12").Value = 0.4427965094191643
$ws.Range("Q12").Value = 1630.964283174493
$ws.Range("R12").Value = 14678.67854857044
$ws.Range("S12").Value = 0.0367599367251279
$ws.Range("T12").Value = 0.03675993672512789
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 79.06597733333334
$ws.Range("H13").Value = 237.197932
$ws.Range("I13").Value = 0.08301767503395074
$ws.Range("J13").Value = 0.08301767503395074
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.125502333333333
$ws.Range("N13").Value = 21.376507
$ws.Range("O13").Value = 0.1529554191465104
$ws.Range("P13").Value = 0.1529554191465104
$ws.Range("Q13").Value = 563.3848059759472
$ws.Range("R13").Value = 5070.463253783524
$ws.Range("S13").Value = 0.01269800328138673
$ws.Range("T13").Value = 0.01269800328138673
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 151.5583673333333
$ws.Range("H14").Value = 454.675102
$ws.Range("I14").Value = 0.1591332164896969
$ws.Range("J14").Value = 0.1591332164896969
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 16.3542
$ws.Range("N14").Value = 49.0626
$ws.Range("O14").Value = 0.3510578481048182
$ws.Range("P14").Value = 0.3510578481048182
$ws.Range("Q14").Value = 2478.6158510428
$ws.Range("R14").Value = 22307.5426593852
$ws.Range("S14").Value = 0.05586496454287117
$ws.Range("T14").Value = 0.05586496454287117
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 151.5583673333333
$ws.Range("H15").Value = 454.675102
$ws.Range("I15").Value = 0.1591332164896969
$ws.Range("J15").Value = 0.1591332164896969
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.477892333333333
$ws.Range("N15").Value = 7.433676999999999
$ws.Range("O15").Value = 0.05319022332950721
$ws.Range("P15").Value = 0.0531902233295072
$ws.Range("Q15").Value = 375.5453164677838
$ws.Range("R15").Value = 3379.907848210054
$ws.Range("S15").Value = 0.008464331324229798
$ws.Range("T15").Value = 0.008464331324229796
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 151.5583673333333
$ws.Range("H16").Value = 454.675102
$ws.Range("I16").Value = 0.1591332164896969
$ws.Range("J16").Value = 0.1591332164896969
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 20.62789
$ws.Range("N16").Value = 61.88367
$ws.Range("O16").Value = 0.4427965094191643
$ws.Range("P16").Value = 0.4427965094191643
$ws.Range("Q16").Value = 3126.329329931593
$ws.Range("R16").Value = 28136.96396938434
$ws.Range("S16").Value = 0.070463632794282
$ws.Range("T16").Value = 0.07046363279428199
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 151.5583673333333
$ws.Range("H17").Value = 454.675102
$ws.Range("I17").Value = 0.1591332164896969
$ws.Range("J17").Value = 0.1591332164896969
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 7.125502333333333
$ws.Range("N17").Value = 21.376507
$ws.Range("O17").Value = 0.1529554191465104
$ws.Range("P17").Value = 0.1529554191465104
$ws.Range("Q17").Value = 1079.929500069857
$ws.Range("R17").Value = 9719.365500628715
$ws.Range("S17").Value = 0.02434028782831398
$ws.Range("T17").Value = 0.02434028782831398
